# Fruta / hortaliza, semanal
# Insert a new weekly record at row 16 (pushing existing rows 16-35 down to 17-36)
# and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 16, shifting rows 16-35 down to 17-36.
$ws.Rows.Item(16).Insert()

# Fill in the new row 16 with this week's record.
$ws.Cells.Item(16, 1).Value = 5
$ws.Cells.Item(16, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(16, 3).Value = "Maule"
$ws.Cells.Item(16, 4).Value = 44763
$ws.Cells.Item(16, 5).Value = 7
$ws.Cells.Item(16, 6).Value = 100112040
$ws.Cells.Item(16, 7).Value = "Cilantro"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 300
$ws.Cells.Item(16, 11).Value = 15000
$ws.Cells.Item(16, 12).Value = 15000
$ws.Cells.Item(16, 13).Value = 15000
$ws.Cells.Item(16, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(16, 15).Value = "Región Metropolitana"
$ws.Cells.Item(16, 16).Value = 417
$ws.Cells.Item(16, 17).Value = 36
$ws.Cells.Item(16, 18).Value = "Hortaliza"
